$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Razon social" entries: a stray comma used inside the name was meant to
#    be a period. These are plain text cells, so a simple in-place Replace
#    is enough (it does not re-infer the cell as a number and does not
#    touch any styles).
# ---------------------------------------------------------------------------
$ws.Range("E35").Replace(",", ".") | Out-Null   # FERNANDEZ MARIO H, GALLICET OSCAR M
$ws.Range("E36").Replace(",", ".") | Out-Null   # IZAGUIRRE CARLOS MARIA, MOREND ...
$ws.Range("E55").Replace(",", ".") | Out-Null   # GIMENEZ, ROBERTO ADRIAN

# ---------------------------------------------------------------------------
# 2) "Importe" column (H2:H76): these amounts were scraped with a
#    comma-decimal / dot-thousands formatting (e.g. "2.725,00"). The fix
#    re-writes them using a plain dot-decimal formatting (e.g. "2725.00"),
#    while keeping the cells as TEXT (they are not real numeric cells in
#    the original workbook).
#
#    Because the resulting text (e.g. "2725.00") looks exactly like a
#    number, a normal assignment would make Excel silently convert the
#    cell to a numeric value (losing the trailing zeros / exact text).
#    Prefixing the value with a leading apostrophe forces Excel to keep it
#    as text; resetting the cell Style back to "Normal" afterwards removes
#    the quote-prefix indicator so the cell ends up with the same (default)
#    style it started with.
# ---------------------------------------------------------------------------
$importes = @(
  "2725.00", "500.00", "1267.82", "1070.00", "2090.00", "75950.00", "172590.19",
  "12652.50", "19799.94", "5613.75", "99.40", "8277.11", "8492.65", "1819.00",
  "6342.16", "929.04", "1195.00", "18566.40", "624.00", "4040.00", "3600.00",
  "162.35", "70.00", "2620.00", "340.00", "273.73", "44975.00", "37128.00",
  "160.00", "1192.50", "1192.50", "15432.00", "952.00", "1638.00", "3165.52",
  "206.00", "50.00", "99951.77", "179.45", "1200.00", "623.29", "10000.00",
  "1862.00", "179.60", "38.40", "3352.00", "2200.00", "1525.93", "330.50",
  "3200.00", "250.00", "1000.00", "420.00", "800.00", "950.00", "2540.00",
  "400.00", "750.00", "200.00", "714.00", "32.30", "7392.00", "706.82",
  "274.00", "8334.20", "4065.00", "2892.62", "849.42", "68750.00",
  "136468.80", "321346.05", "17000.00", "33300.00", "40000.00"
)

$row = 2
foreach ($val in $importes) {
  $cell = $ws.Cells.Item($row, 8)   # column H
  $cell.Value = "'" + $val
  $cell.Style = "Normal"
  $row++
}
